# edit.ps1 - reproduces the author's edit:
#   1. Fixes the GitHub profile URL casing on the last slide
#      ("Surabhi0901" -> "surabhi0901") while leaving the rest of the
#      hyperlink run's formatting (hyperlink target, underline, size)
#      untouched.
#   2. Restores the presentation's applied color theme to the original
#      "Default" palette (the palette that shipped in the deck's other,
#      currently-unused theme part) instead of the "Material" palette
#      that is presently applied to the slide master/deck.

$p = $ppt.ActivePresentation

# --- 1. Fix the GitHub link casing on the last slide ---------------------
$lastSlideIndex = $p.Slides.Count
$slide = $p.Slides.Item($lastSlideIndex)

# The "Contact me" body placeholder is the second shape on this slide.
$contactShape = $slide.Shapes.Item(2)
$tr = $contactShape.TextFrame.TextRange

$oldUrl = "https://github.com/Surabhi0901"
$newUrl = "https://github.com/surabhi0901"

$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldUrl) + 1
if ($startPos -gt 0) {
    $target = $tr.Characters($startPos, $oldUrl.Length)
    $target.Text = $newUrl
}

# --- 2. Revert the deck's applied theme colors back to the original ------
# "Default" scheme (currently the "Material" scheme is applied). The
# font/format schemes are identical between the two themes already, so
# only the 12 theme colors need to change.
$themeColors = $slide.ThemeColorScheme

# index -> (name, RGB as 0xBBGGRR packed int, i.e. VBA/COM RGB())
$themeColors.Item(1).RGB  = 0          # dk1      000000
$themeColors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Item(3).RGB  = 5800213    # dk2      158158
$themeColors.Item(4).RGB  = 15987699   # lt2      F3F3F3
$themeColors.Item(5).RGB  = 13077765   # accent1  058DC7
$themeColors.Item(6).RGB  = 3322960    # accent2  50B432
$themeColors.Item(7).RGB  = 1791725    # accent3  ED561B
$themeColors.Item(8).RGB  = 61421      # accent4  EDEF00
$themeColors.Item(9).RGB  = 15059748   # accent5  24CBE5
$themeColors.Item(10).RGB = 7529828    # accent6  64E572
$themeColors.Item(11).RGB = 13369378   # hlink    2200CC
$themeColors.Item(12).RGB = 9116245    # folHlink 551A8B

Write-Output ("Updated contact text: " + $tr.Text)
Write-Output ("Theme color 1 (dk1) now: " + $themeColors.Item(1).RGB)
